# Update the "想去人数" (number of people interested) counts that changed
# between data snapshots for two events, in both the "展览" and "全部类型"
# worksheets (which mirror the same rows of data).
#   F4: 60   -> 64
#   F5: 2459 -> 2476

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 64
    $ws.Range("F5").Value = 2476
}
